$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 131, shifting existing rows 131-211 down to 133-213
$ws.Rows('131:132').Insert()

$rowData = @{
  131 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44582, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 800, 900, 850, '$/unidad', 'Región Metropolitana', 850, 1, 'Hortaliza')
  132 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44582, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 700, 700, 700, '$/unidad', 'Región Metropolitana', 700, 1, 'Hortaliza')
  133 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44308, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  134 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44308, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  135 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44357, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  136 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44357, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  137 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44320, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  138 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44320, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  139 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44306, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1500, 600, 700, 647, '$/unidad', 'Región Metropolitana', 647, 1, 'Hortaliza')
  140 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44306, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 600, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  141 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44295, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  142 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44295, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
  143 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44210, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 600, 700, 650, '$/unidad', 'Región del Maule', 650, 1, 'Hortaliza')
  144 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44210, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 500, 500, 500, '$/unidad', 'Región del Maule', 500, 1, 'Hortaliza')
  145 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44343, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  146 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44343, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  147 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44316, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  148 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44316, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  149 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44329, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  150 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44329, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  151 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44460, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  152 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44460, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  153 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44526, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  154 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44526, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
  155 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44313, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  156 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44313, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 800, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  157 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44334, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  158 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44334, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  159 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44483, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2100, 600, 650, 629, '$/unidad', 'Provincia de Melipilla', 629, 1, 'Hortaliza')
  160 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44448, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  161 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44448, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
  162 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44196, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región del Maule', 750, 1, 'Hortaliza')
  163 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44196, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 600, 600, 600, '$/unidad', 'Región del Maule', 600, 1, 'Hortaliza')
  164 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44463, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  165 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44463, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
  166 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44476, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  167 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44476, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  168 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44496, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1050, 750, 750, 750, '$/unidad', 'Región del Maule', 750, 1, 'Hortaliza')
  169 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44208, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  170 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44208, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 300, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
  171 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44250, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 900, 1000, 950, '$/unidad', 'Región Metropolitana', 950, 1, 'Hortaliza')
  172 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44250, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 800, 800, 800, '$/unidad', 'Región Metropolitana', 800, 1, 'Hortaliza')
  173 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44560, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  174 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44560, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  175 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44216, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 600, 700, 650, '$/unidad', 'Región del Maule', 650, 1, 'Hortaliza')
  176 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44216, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 500, 500, 500, '$/unidad', 'Región del Maule', 500, 1, 'Hortaliza')
  177 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44509, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  178 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44509, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
  179 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44229, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  180 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44229, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  181 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44488, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  182 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44488, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
  183 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44341, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 650, 625, '$/unidad', 'Región Metropolitana', 625, 1, 'Hortaliza')
  184 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44341, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  185 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44482, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 850, 800, 900, 853, '$/unidad', 'Provincia del Elquí', 853, 1, 'Hortaliza')
  186 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44558, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1250, 700, 750, 726, '$/unidad', 'Región Metropolitana', 726, 1, 'Hortaliza')
  187 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44558, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 650, 500, 550, 523, '$/unidad', 'Región Metropolitana', 523, 1, 'Hortaliza')
  188 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44545, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2700, 400, 450, 428, '$/unidad', 'Región Metropolitana', 428, 1, 'Hortaliza')
  189 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44322, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 700, 800, 750, '$/unidad', 'Región del Maule', 750, 1, 'Hortaliza')
  190 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44322, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 600, 600, 600, '$/unidad', 'Región del Maule', 600, 1, 'Hortaliza')
  191 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44194, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 800, 900, 850, '$/unidad', 'Región Metropolitana', 850, 1, 'Hortaliza')
  192 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44194, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 700, 700, 700, '$/unidad', 'Región Metropolitana', 700, 1, 'Hortaliza')
  193 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44236, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región del Maule', 750, 1, 'Hortaliza')
  194 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44236, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 600, 600, 600, '$/unidad', 'Región del Maule', 600, 1, 'Hortaliza')
  195 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44162, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 800, 600, 700, 650, '$/unidad', 'Región del Maule', 650, 1, 'Hortaliza')
  196 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44162, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 400, 500, 500, 500, '$/unidad', 'Región del Maule', 500, 1, 'Hortaliza')
  197 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44427, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  198 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44427, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  199 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44491, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 800, 900, 850, '$/unidad', 'Región Metropolitana', 850, 1, 'Hortaliza')
  200 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44491, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 700, 700, 700, '$/unidad', 'Región Metropolitana', 700, 1, 'Hortaliza')
  201 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44266, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 800, 700, 750, 725, '$/unidad', 'Región del Maule', 725, 1, 'Hortaliza')
  202 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44266, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 400, 600, 600, 600, '$/unidad', 'Región del Maule', 600, 1, 'Hortaliza')
  203 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44533, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1500, 600, 650, 623, '$/unidad', 'Región Metropolitana', 623, 1, 'Hortaliza')
  204 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44264, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 800, 900, 850, '$/unidad', 'Región Metropolitana', 850, 1, 'Hortaliza')
  205 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44264, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 700, 700, 700, '$/unidad', 'Región Metropolitana', 700, 1, 'Hortaliza')
  206 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44279, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 700, 800, 750, '$/unidad', 'Región del Maule', 750, 1, 'Hortaliza')
  207 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44279, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 600, 600, 600, '$/unidad', 'Región del Maule', 600, 1, 'Hortaliza')
  208 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44327, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  209 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44327, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  210 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44462, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 1000, 600, 700, 650, '$/unidad', 'Región Metropolitana', 650, 1, 'Hortaliza')
  211 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44462, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 500, 500, 500, 500, '$/unidad', 'Región Metropolitana', 500, 1, 'Hortaliza')
  212 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44512, 8, 100112023, 'Brócoli', 'Sin especificar', 'Primera', 2000, 700, 800, 750, '$/unidad', 'Región Metropolitana', 750, 1, 'Hortaliza')
  213 = @(11, 'Vega Monumental Concepción', 'Bíobío', 44512, 8, 100112023, 'Brócoli', 'Sin especificar', 'Segunda', 1000, 600, 600, 600, '$/unidad', 'Región Metropolitana', 600, 1, 'Hortaliza')
}

foreach ($rowNum in $rowData.Keys) {
  $vals = $rowData[$rowNum]
  for ($c = 0; $c -lt $vals.Count; $c++) {
    $ws.Cells.Item([int]$rowNum, $c + 1).Value = $vals[$c]
  }
}
